$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.269101333333333
$ws.Range("H2").Value = 6.807304
$ws.Range("I2").Value = 0.02891211995713196
$ws.Range("J2").Value = 0.02891211995713196
$ws.Range("M2").Value = 28.689524
$ws.Range("N2").Value = 86.06857199999999
$ws.Range("O2").Value = 0.2394085694101769
$ws.Range("P2").Value = 0.2394085694101769
$ws.Range("Q2").Value = 65.09943716109866
$ws.Range("R2").Value = 585.894934449888
$ws.Range("S2").Value = 0.006921809277552388
$ws.Range("T2").Value = 0.006921809277552388
$ws.Range("G3").Value = 2.269101333333333
$ws.Range("H3").Value = 6.807304
$ws.Range("I3").Value = 0.02891211995713196
$ws.Range("J3").Value = 0.02891211995713196
$ws.Range("O3").Value = 0.5212694246546397
$ws.Range("P3").Value = 0.5212694246546395
$ws.Range("Q3").Value = 141.7424039494898
$ws.Range("R3").Value = 1275.681635545408
$ws.Range("S3").Value = 0.0150710041356001
$ws.Range("T3").Value = 0.0150710041356001
$ws.Range("G4").Value = 2.269101333333333
$ws.Range("H4").Value = 6.807304
$ws.Range("I4").Value = 0.02891211995713196
$ws.Range("J4").Value = 0.02891211995713196
$ws.Range("M4").Value = 28.525746
$ws.Range("N4").Value = 85.57723799999999
$ws.Range("O4").Value = 0.2380418734454457
$ws.Range("P4").Value = 0.2380418734454457
$ws.Range("Q4").Value = 64.727808282928
$ws.Range("R4").Value = 582.5502745463519
$ws.Range("S4").Value = 0.006882295199875151
$ws.Range("T4").Value = 0.00688229519987515
$ws.Range("G5").Value = 2.269101333333333
$ws.Range("H5").Value = 6.807304
$ws.Range("I5").Value = 0.02891211995713196
$ws.Range("J5").Value = 0.02891211995713196
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1534046666666667
$ws.Range("N5").Value = 0.460214
$ws.Range("O5").Value = 0.001280132489737778
$ws.Range("P5").Value = 0.001280132489737778
$ws.Range("Q5").Value = 0.3480907336728889
$ws.Range("R5").Value = 3.132816603056
$ws.Range("S5").Value = 0.00003701134410432062
$ws.Range("T5").Value = 0.00003701134410432062
$ws.Range("I6").Value = 0.7238963226334669
$ws.Range("J6").Value = 0.7238963226334669
$ws.Range("M6").Value = 28.689524
$ws.Range("N6").Value = 86.06857199999999
$ws.Range("O6").Value = 0.2394085694101769
$ws.Range("P6").Value = 0.2394085694101769
$ws.Range("Q6").Value = 1629.947691013334
$ws.Range("R6").Value = 14669.52921912001
$ws.Range("S6").Value = 0.1733069830029662
$ws.Range("T6").Value = 0.1733069830029662
$ws.Range("I7").Value = 0.7238963226334669
$ws.Range("J7").Value = 0.7238963226334669
$ws.Range("O7").Value = 0.5212694246546397
$ws.Range("P7").Value = 0.5212694246546395
$ws.Range("S7").Value = 0.3773450196087567
$ws.Range("T7").Value = 0.3773450196087566
$ws.Range("I8").Value = 0.7238963226334669
$ws.Range("J8").Value = 0.7238963226334669
$ws.Range("M8").Value = 28.525746
$ws.Range("N8").Value = 85.57723799999999
$ws.Range("O8").Value = 0.2380418734454457
$ws.Range("P8").Value = 0.2380418734454457
$ws.Range("Q8").Value = 1620.642915760222
$ws.Range("R8").Value = 14585.786241842
$ws.Range("S8").Value = 0.1723176368199393
$ws.Range("T8").Value = 0.1723176368199393
$ws.Range("I9").Value = 0.7238963226334669
$ws.Range("J9").Value = 0.7238963226334669
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1534046666666667
$ws.Range("N9").Value = 0.460214
$ws.Range("O9").Value = 0.001280132489737778
$ws.Range("P9").Value = 0.001280132489737778
$ws.Range("Q9").Value = 8.715431536054888
$ws.Range("R9").Value = 78.438883824494
$ws.Range("S9").Value = 0.0009266832018048016
$ws.Range("T9").Value = 0.0009266832018048016
$ws.Range("G10").Value = 18.57257166666666
$ws.Range("H10").Value = 55.717715
$ws.Range("I10").Value = 0.2366454120188096
$ws.Range("J10").Value = 0.2366454120188096
$ws.Range("M10").Value = 28.689524
$ws.Range("N10").Value = 86.06857199999999
$ws.Range("O10").Value = 0.2394085694101769
$ws.Range("P10").Value = 0.2394085694101769
$ws.Range("Q10").Value = 532.8382405725532
$ws.Range("R10").Value = 4795.544165152979
$ws.Range("S10").Value = 0.05665493954890509
$ws.Range("T10").Value = 0.0566549395489051
$ws.Range("G11").Value = 18.57257166666666
$ws.Range("H11").Value = 55.717715
$ws.Range("I11").Value = 0.2366454120188096
$ws.Range("J11").Value = 0.2366454120188096
$ws.Range("O11").Value = 0.5212694246546397
$ws.Range("P11").Value = 0.5212694246546395
$ws.Range("Q11").Value = 1160.160155426076
$ws.Range("R11").Value = 10441.44139883468
$ws.Range("S11").Value = 0.123356017770205
$ws.Range("T11").Value = 0.123356017770205
$ws.Range("G12").Value = 18.57257166666666
$ws.Range("H12").Value = 55.717715
$ws.Range("I12").Value = 0.2366454120188096
$ws.Range("J12").Value = 0.2366454120188096
$ws.Range("M12").Value = 28.525746
$ws.Range("N12").Value = 85.57723799999999
$ws.Range("O12").Value = 0.2380418734454457
$ws.Range("P12").Value = 0.2380418734454457
$ws.Range("Q12").Value = 529.7964619301299
$ws.Range("R12").Value = 4768.168157371169
$ws.Range("S12").Value = 0.05633151721922683
$ws.Range("T12").Value = 0.05633151721922683
$ws.Range("G13").Value = 18.57257166666666
$ws.Range("H13").Value = 55.717715
$ws.Range("I13").Value = 0.2366454120188096
$ws.Range("J13").Value = 0.2366454120188096
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1534046666666667
$ws.Range("N13").Value = 0.460214
$ws.Range("O13").Value = 0.001280132489737778
$ws.Range("P13").Value = 0.001280132489737778
$ws.Range("Q13").Value = 2.849119165667777
$ws.Range("R13").Value = 25.64207249101
$ws.Range("S13").Value = 0.0003029374804726609
$ws.Range("T13").Value = 0.0003029374804726609
$ws.Range("G14").Value = 0.8276899999999999
$ws.Range("H14").Value = 2.48307
$ws.Range("I14").Value = 0.01054614539059158
$ws.Range("J14").Value = 0.01054614539059158
$ws.Range("M14").Value = 28.689524
$ws.Range("N14").Value = 86.06857199999999
$ws.Range("O14").Value = 0.2394085694101769
$ws.Range("P14").Value = 0.2394085694101769
$ws.Range("Q14").Value = 23.74603211956
$ws.Range("R14").Value = 213.7142890760399
$ws.Range("S14").Value = 0.002524837580753262
$ws.Range("T14").Value = 0.002524837580753262
$ws.Range("G15").Value = 0.8276899999999999
$ws.Range("H15").Value = 2.48307
$ws.Range("I15").Value = 0.01054614539059158
$ws.Range("J15").Value = 0.01054614539059158
$ws.Range("O15").Value = 0.5212694246546397
$ws.Range("P15").Value = 0.5212694246546395
$ws.Range("Q15").Value = 51.70274619362667
$ws.Range("R15").Value = 465.32471574264
$ws.Range("S15").Value = 0.005497383140077854
$ws.Range("T15").Value = 0.005497383140077853
$ws.Range("G16").Value = 0.8276899999999999
$ws.Range("H16").Value = 2.48307
$ws.Range("I16").Value = 0.01054614539059158
$ws.Range("J16").Value = 0.01054614539059158
$ws.Range("M16").Value = 28.525746
$ws.Range("N16").Value = 85.57723799999999
$ws.Range("O16").Value = 0.2380418734454457
$ws.Range("P16").Value = 0.2380418734454457
$ws.Range("Q16").Value = 23.61047470674
$ws.Range("R16").Value = 212.4942723606599
$ws.Range("S16").Value = 0.002510424206404472
$ws.Range("T16").Value = 0.002510424206404472
$ws.Range("G17").Value = 0.8276899999999999
$ws.Range("H17").Value = 2.48307
$ws.Range("I17").Value = 0.01054614539059158
$ws.Range("J17").Value = 0.01054614539059158
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1534046666666667
$ws.Range("N17").Value = 0.460214
$ws.Range("O17").Value = 0.001280132489737778
$ws.Range("P17").Value = 0.001280132489737778
$ws.Range("Q17").Value = 0.1269715085533333
$ws.Range("R17").Value = 1.14274357698
$ws.Range("S17").Value = 0.00001350046335599459
$ws.Range("T17").Value = 0.00001350046335599459
